# Weekly update: insert a new data row for the latest date at the top of
# the data block (row 4), pushing all existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 4 (the first data row after the
# two most-recent existing rows). This shifts rows 4..112 down to 5..113
# and keeps formatting consistent with the rest of the table (Excel copies
# the formatting of the row above, which already uses the date style
# needed for column D).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the latest record.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 45160
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107011
$ws.Range("J4").Value = "Tuna"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 33000
$ws.Range("O4").Value = 33000
$ws.Range("P4").Value = 33000
$ws.Range("Q4").Value = "$/caja 16 kilos"
$ws.Range("R4").Value = "Provincia de Los Andes"
$ws.Range("S4").Value = 2062
$ws.Range("T4").Value = 16

# Make sure the date cell keeps/has the correct date formatting.
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
